$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.220278859138489
$ws.Range("B1").Value = 1.930655121803284
$ws.Range("C1").Value = 4.239514350891113
$ws.Range("D1").Value = 3.101081848144531
$ws.Range("E1").Value = 1.187941551208496
